$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 needs to hold the text "4684" (not a number). A plain .Value assignment
# of a numeric-looking string gets auto-coerced to a number, so stage it in
# a scratch cell formatted as Text, copy/paste-special the value across,
# then wipe the scratch cell (including its formatting) so nothing extra
# is left behind in the workbook.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "4684"
$ws.Range("Z1").Copy()
$ws.Range("A2").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("C2").Value = "bat"
$ws.Range("G2").Value = "29S"

$ws.Range("B2").Value = 14.0
$ws.Range("F2").Value = 1460.66797549867
$ws.Range("H2").Value = 0.03166666666666666
$ws.Range("I2").Value = 2.0316666666666667
